$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new row of data (row 11) ---
$ws.Cells.Item(11, 1).Value2 = "SER25-00777"
$ws.Cells.Item(11, 2).Value2 = 60
$ws.Cells.Item(11, 3).Value2 = "mermaids"
$ws.Cells.Item(11, 4).Value2 = "Tursiops"
$ws.Cells.Item(11, 5).Value2 = "bellissima"
$ws.Cells.Item(11, 6).Value2 = "test 4"

# Copy the formatting (border / fill / font) from row 10 down into row 11,
# matching the bordered "data row" look used throughout the table.
$ws.Range("A10:F10").Copy() | Out-Null
$ws.Range("A11:F11").PasteSpecial(-4122) | Out-Null

# Column E keeps the plain (unbordered) style used elsewhere in the sheet,
# so restore it from E9 after the row-wide paste above.
$ws.Range("E9").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Row heights nudge up slightly (16.85 -> 16.9) across the table once the new
# row is added.
$ws.Range("A1:F11").RowHeight = 16.9

# --- Selection moves to C15 ---
$ws.Range("C15").Select() | Out-Null
